# Apply edits described by the diff:
#  - TestSuite sheet: change RunMode for AdminPageTest (B3) from "N" to "Y"
#                     and update the active cell selection to B3
#  - LoginPageTest sheet: update the active cell selection to C3

$wb = $excel.ActiveWorkbook

# --- TestSuite sheet ---
$wsSuite = $wb.Worksheets.Item("TestSuite")
$wsSuite.Activate()
$wsSuite.Range("B3").Value = "Y"
$wsSuite.Range("B3").Select()

# --- LoginPageTest sheet ---
$wsLogin = $wb.Worksheets.Item("LoginPageTest")
$wsLogin.Activate()
$wsLogin.Range("C3").Select()

# Re-activate TestSuite sheet as the final active/selected sheet (tabSelected="1")
$wsSuite.Activate()
